$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename two header cells to reflect the "generic rename" work:
# M2 "tested" -> "finished and tested"
# H2 "Completion method written" -> "Completion or callback written"
# (Set M2 first so the new shared strings are appended in the same order
# as the target workbook: "finished and tested" then "Completion or
# callback written".)
$ws.Range("M2").Value = "finished and tested"
$ws.Range("H2").Value = "Completion or callback written"

# Fill in the previously-blank K/L columns ("js modified" / "jade modifed")
# for the rows that already use the Generic Rename logic (n/a) ...
$ws.Range("K3").Value = "n/a"
$ws.Range("L3").Value = "n/a"

$ws.Range("K6").Value = "n/a"
$ws.Range("L6").Value = "n/a"

$ws.Range("K7").Value = "n/a"
$ws.Range("L7").Value = "n/a"

$ws.Range("K11").Value = "n/a"
$ws.Range("L11").Value = "n/a"

$ws.Range("K12").Value = "n/a"
$ws.Range("L12").Value = "n/a"

$ws.Range("K14").Value = "n/a"
$ws.Range("L14").Value = "n/a"

$ws.Range("K15").Value = "n/a"
$ws.Range("L15").Value = "n/a"

# ... and "x" for the confirmation / generic rename definition rows.
$ws.Range("K16").Value = "x"
$ws.Range("L16").Value = "x"

$ws.Range("K17").Value = "x"
$ws.Range("L17").Value = "x"

# Move the active selection from K7 to H3.
$ws.Range("H3").Select()
